# "Generate Report for Handback" - mark the two handed-back files as
# synced and point the "Latest Target File" column at the markdown that
# was handed back, with its handback xlf name + timestamp.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdUrl5560 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11050a68bce6cf96d392cd8929cac344cedaea58/e2e/5560ddf0-dae4-4751-aaf6-75e9136f135a.md"
$mdUrlFe70 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11050a68bce6cf96d392cd8929cac344cedaea58/e2e/fe708d45-9349-4462-a8b3-91111eff3bad.md"
$md5560 = "5560ddf0-dae4-4751-aaf6-75e9136f135a.md"
$mdFe70 = "fe708d45-9349-4462-a8b3-91111eff3bad.md"

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status columns (E, F) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl5560, "", "", $md5560)
$wsZh.Range("J2").Value = "5560ddf0-dae4-4751-aaf6-75e9136f135a.b7a04fe5db1ebb8cede6d09ac74d860e1a066298.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 16:35:56"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrlFe70, "", "", $mdFe70)
$wsZh.Range("J3").Value = "fe708d45-9349-4462-a8b3-91111eff3bad.52e29b1a62bc63afacb842f0b94461756ed8ef40.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-30 16:35:56"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl5560, "", "", $md5560)
$wsDe.Range("J2").Value = "5560ddf0-dae4-4751-aaf6-75e9136f135a.b7a04fe5db1ebb8cede6d09ac74d860e1a066298.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 16:36:15"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrlFe70, "", "", $mdFe70)
$wsDe.Range("J3").Value = "fe708d45-9349-4462-a8b3-91111eff3bad.52e29b1a62bc63afacb842f0b94461756ed8ef40.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-30 16:36:15"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
